# "Actualizar" (availability-check) run appended on 2021-02-05 ~10:18:52.
# Mirrors the existing 14-row cycle (Odoo/Blackbox/PowerBI/Dropbox/Odoo/GEE/
# UtilidadesOdoo/Filtros Dashboard/MapStore/GeoServer/Tomcat/Shiny/Github/
# EZ Exporter) that every previous "Actualizar" batch wrote to this sheet,
# plus a small timestamp correction on the immediately preceding batch
# (rows 492-505) that the same commit carried.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fix the timestamp on the previous batch (rows 492-505) ---------
for ($r = 492; $r -le 505; $r++) {
    $ws.Cells($r, 4).Value = 44232.40867990741
}

# --- 2. Append the new batch: rows 506-519 ------------------------------
$names = @("Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", `
           "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer", `
           "Tomcat", "Shiny", "Github", "EZ Exporter")

# Address actually used for the relationship target (no fragment for MapStore)
$addresses = @( `
    "https://www.dataintelligence-group.com/", `
    "https://serviciodashboard.azurewebsites.net/", `
    "https://powerbi.microsoft.com/es-es/", `
    "https://www.dropbox.com/", `
    "https://dataintelligence.store/", `
    "https://app-data-i.users.earthengine.app/", `
    "https://odooutil.azurewebsites.net/", `
    "https://filtradordashboard.azurewebsites.net/", `
    "https://ide.dataintelligence-group.com/mapstore/", `
    "https://ide.dataintelligence-group.com/geoserver/web/?0", `
    "https://ide.dataintelligence-group.com/", `
    "https://rpubs.com/dataintelligence/", `
    "https://github.com/Sud-Austral/", `
    "https://ezexporter.highviewapps.com/exports/export-profile/" `
)

# Sub-address (in-page fragment) tacked on the displayed/stored text; only
# the MapStore row (index 8) carries one ("#/" -> stored as location "/").
$subAddresses = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")

$newTimestamp = 44232.42977864746

for ($i = 0; $i -lt 14; $i++) {
    $row = 506 + $i

    $displayUrl = $addresses[$i]
    if ($subAddresses[$i] -ne "") {
        $displayUrl = $displayUrl + "#" + $subAddresses[$i]
    }

    $ws.Cells($row, 1).Value = $names[$i]
    $ws.Cells($row, 2).Value = $displayUrl
    $ws.Cells($row, 3).Value = "Disponible"
    $ws.Cells($row, 4).Value = $newTimestamp

    if ($subAddresses[$i] -ne "") {
        $ws.Hyperlinks.Add($ws.Cells($row, 2), $addresses[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($ws.Cells($row, 2), $addresses[$i])
    }

    # Hyperlinks.Add stamps its own ad-hoc font xf; put column B back on the
    # workbook's existing shared "Hyperlink" cell style (same xf every other
    # link cell in the sheet already uses).
    $ws.Cells($row, 2).Style = "Hyperlink"

    # Column D keeps the workbook's custom date/time number format.
    $ws.Cells($row, 4).NumberFormat = $ws.Cells($row - 1, 4).NumberFormat
}

Write-Output "Applied Actualizar batch (rows 506-519) + timestamp fix (492-505)"
